$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 1.336873824401267

$ws.Range("B3").Value = 0.3464964993005633
$ws.Range("C3").Value = 0.05231270169004087
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 7.592970129276214

$ws.Range("B4").Value = 0.7287194209349384
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 3.034748368925986

$ws.Range("B5").Value = 0.7287194209349384
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 3.082599426703578
$ws.Range("E5").Value = 6.48142807727062
$ws.Range("G5").Value = 11.94598338380795

$ws.Range("B6").Value = 3.182878228561681
$ws.Range("C6").Value = 1.65323645889881
$ws.Range("D6").Value = 0.7127328510149897
$ws.Range("E6").Value = 0.4998867070740569
$ws.Range("G6").Value = 6.048734245549538

$ws.Range("B7").Value = 0.7287194209349384
$ws.Range("C7").Value = 1.65323645889881
$ws.Range("D7").Value = 0.1529057820181812
$ws.Range("E7").Value = 0.4998867070740569
$ws.Range("G7").Value = 3.034748368925986

$ws.Range("B8").Value = 1.505614041169197
$ws.Range("C8").Value = 1.65323645889881
$ws.Range("D8").Value = 0.1529057820181812
$ws.Range("E8").Value = 6.48142807727062
$ws.Range("G8").Value = 9.793184359356808

$ws.Range("B9").Value = 0.1554434735375247
$ws.Range("C9").Value = 0.0001537489499301437
$ws.Range("D9").Value = 0.1529057820181812
$ws.Range("E9").Value = 0.4998867070740569
$ws.Range("G9").Value = 0.808389711579693

$ws.Range("B10").Value = 0.3464964993005633
$ws.Range("C10").Value = 1.65323645889881
$ws.Range("D10").Value = 0.7127328510149897
$ws.Range("E10").Value = 6.48142807727062
$ws.Range("G10").Value = 9.193893886484982

$ws.Range("B11").Value = 1.505614041169197
$ws.Range("C11").Value = 1.65323645889881
$ws.Range("D11").Value = 0.1529057820181812
$ws.Range("E11").Value = 0.4998867070740569
$ws.Range("G11").Value = 3.811642989160245

$ws.Range("B12").Value = 0.0003714022599530242
$ws.Range("C12").Value = 0.0001537489499301437
$ws.Range("D12").Value = 0.1529057820181812
$ws.Range("E12").Value = 0.4998867070740569
$ws.Range("G12").Value = 0.6533176403021212
